$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1)
$ws.Range('A1').Value = 'mx_state'
$ws.Range('B1').Value = 'mx_municipality'
$ws.Range('C1').Value = 'n_matriculas'
$ws.Range('D1').Value = 'pct_matriculas'

# Convert state/municipality names from UPPERCASE to Title Case
$ws.Range('A2').Value = 'Aguascalientes'
$ws.Range('B2').Value = 'Aguascalientes'
$ws.Range('B3').Value = 'Pabellón De Arteaga'
$ws.Range('B4').Value = 'Total'
$ws.Range('A5').Value = 'Baja California'
$ws.Range('B5').Value = 'Tijuana'
$ws.Range('B6').Value = 'Total'
$ws.Range('A7').Value = 'Campeche'
$ws.Range('B7').Value = 'Carmen'
$ws.Range('B8').Value = 'Total'
$ws.Range('A9').Value = 'Chiapas'
$ws.Range('B9').Value = 'Amatenango De La Frontera'
$ws.Range('B10').Value = 'Escuintla'
$ws.Range('B11').Value = 'Palenque'
$ws.Range('B12').Value = 'Tapachula'
$ws.Range('B13').Value = 'Tonalá'
$ws.Range('B14').Value = 'Villa Corzo'
$ws.Range('B15').Value = 'Total'
$ws.Range('A16').Value = 'Chihuahua'
$ws.Range('B16').Value = 'Juárez'
$ws.Range('B17').Value = 'Total'
$ws.Range('A18').Value = 'Ciudad De México'
$ws.Range('B18').Value = 'Benito Juárez'
$ws.Range('B19').Value = 'Cuauhtémoc'
$ws.Range('B20').Value = 'Gustavo A. Madero'
$ws.Range('B21').Value = 'Iztapalapa'
$ws.Range('B22').Value = 'La Magdalena Contreras'
$ws.Range('B23').Value = 'Miguel Hidalgo'
$ws.Range('B24').Value = 'No Se Registró El Municipio/Condado/Alcaldía De Nacimiento'
$ws.Range('B25').Value = 'Venustiano Carranza'
$ws.Range('B26').Value = 'Álvaro Obregón'
$ws.Range('B27').Value = 'Total'
$ws.Range('A28').Value = 'Estado De México'
$ws.Range('B28').Value = 'Ixtlahuaca'
$ws.Range('B29').Value = 'Morelos'
$ws.Range('B30').Value = 'Nezahualcóyotl'
$ws.Range('B31').Value = 'San Felipe Del Progreso'
$ws.Range('B32').Value = 'San José Del Rincón'
$ws.Range('B33').Value = 'Temascalcingo'
$ws.Range('B34').Value = 'Teoloyucan'
$ws.Range('B35').Value = 'Toluca'
$ws.Range('B36').Value = 'Tultitlán'
$ws.Range('B37').Value = 'Total'
$ws.Range('A38').Value = 'Guerrero'
$ws.Range('B38').Value = 'Acapulco De Juárez'
$ws.Range('B39').Value = 'Acatepec'
$ws.Range('B40').Value = 'Ahuacuotzingo'
$ws.Range('B41').Value = 'Alcozauca De Guerrero'
$ws.Range('B42').Value = 'Atenango Del Río'
$ws.Range('B43').Value = 'Atlixtac'
$ws.Range('B44').Value = 'Chilapa De Álvarez'
$ws.Range('B45').Value = 'Chilpancingo De Los Bravo'
$ws.Range('B46').Value = 'Coahuayutla De José María Izazaga'
$ws.Range('B47').Value = 'Copanatoyac'
$ws.Range('B48').Value = 'Cualác'
$ws.Range('B49').Value = 'Eduardo Neri'
$ws.Range('B50').Value = 'Olinalá'
$ws.Range('B51').Value = 'Ometepec'
$ws.Range('B52').Value = 'San Miguel Totolapan'
$ws.Range('B53').Value = 'Tlacoapa'
$ws.Range('B54').Value = 'Tlapa De Comonfort'
$ws.Range('B55').Value = 'Zitlala'
$ws.Range('B56').Value = 'Total'
$ws.Range('A57').Value = 'Hidalgo'
$ws.Range('B57').Value = 'Apan'
$ws.Range('B58').Value = 'Total'
$ws.Range('A59').Value = 'Michoacán De Ocampo'
$ws.Range('B59').Value = 'Apatzingán'
$ws.Range('B60').Value = 'Total'
$ws.Range('A61').Value = 'Morelos'
$ws.Range('B61').Value = 'Cuautla'
$ws.Range('B62').Value = 'Tlaquiltenango'
$ws.Range('B63').Value = 'Yecapixtla'
$ws.Range('B64').Value = 'Zacatepec'
$ws.Range('B65').Value = 'Total'
$ws.Range('A66').Value = 'Oaxaca'
$ws.Range('B66').Value = 'Chahuites'
$ws.Range('B67').Value = 'Chalcatongo De Hidalgo'
$ws.Range('B68').Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range('B69').Value = 'Oaxaca De Juárez'
$ws.Range('B70').Value = 'San Francisco Del Mar'
$ws.Range('B71').Value = 'San Francisco Telixtlahuaca'
$ws.Range('B72').Value = 'San Juan Bautista Jayacatlán'
$ws.Range('B73').Value = 'San Juan Bautista Tuxtepec'
$ws.Range('B74').Value = 'San Juan Bautista Valle Nacional'
$ws.Range('B75').Value = 'San Pablo Huitzo'
$ws.Range('B76').Value = 'San Pedro Pochutla'
$ws.Range('B77').Value = 'Santa Catarina Lachatao'
$ws.Range('B78').Value = 'Santiago Jamiltepec'
$ws.Range('B79').Value = 'Total'
$ws.Range('A80').Value = 'Puebla'
$ws.Range('B80').Value = 'Acatzingo'
$ws.Range('B81').Value = 'Chapulco'
$ws.Range('B82').Value = 'Chichiquila'
$ws.Range('B83').Value = 'Epatlán'
$ws.Range('B84').Value = 'General Felipe Ángeles'
$ws.Range('B85').Value = 'Ixcamilpa De Guerrero'
$ws.Range('B86').Value = 'Ixcaquixtla'
$ws.Range('B87').Value = 'Izúcar De Matamoros'
$ws.Range('B88').Value = 'Los Reyes De Juárez'
$ws.Range('B89').Value = 'Nealtican'
$ws.Range('B90').Value = 'Santiago Miahuatlán'
$ws.Range('B91').Value = 'Tehuacán'
$ws.Range('B92').Value = 'Tehuitzingo'
$ws.Range('B93').Value = 'Tepanco De López'
$ws.Range('B94').Value = 'Tlacotepec De Benito Juárez'
$ws.Range('B95').Value = 'Vicente Guerrero'
$ws.Range('B96').Value = 'Zacapoaxtla'
$ws.Range('B97').Value = 'Total'
$ws.Range('A98').Value = 'Querétaro'
$ws.Range('B98').Value = 'San Juan Del Río'
$ws.Range('B99').Value = 'Total'
$ws.Range('A100').Value = 'Quintana Roo'
$ws.Range('B100').Value = 'Othón P. Blanco'
$ws.Range('B101').Value = 'Total'
$ws.Range('A102').Value = 'Sonora'
$ws.Range('B102').Value = 'Alamos'
$ws.Range('B103').Value = 'Total'
$ws.Range('A104').Value = 'Tamaulipas'
$ws.Range('B104').Value = 'Tampico'
$ws.Range('B105').Value = 'Total'
$ws.Range('A106').Value = 'Tlaxcala'
$ws.Range('B106').Value = 'Apizaco'
$ws.Range('B107').Value = 'Tepeyanco'
$ws.Range('B108').Value = 'Tetla De La Solidaridad'
$ws.Range('B109').Value = 'Tlaxcala'
$ws.Range('B110').Value = 'Tocatlán'
$ws.Range('B111').Value = 'Yauhquemehcan'
$ws.Range('B112').Value = 'Total'
$ws.Range('A113').Value = 'Veracruz De Ignacio De La Llave'
$ws.Range('B113').Value = 'Amatlán De Los Reyes'
$ws.Range('B114').Value = 'Córdoba'
$ws.Range('B115').Value = 'Martínez De La Torre'
$ws.Range('B116').Value = 'San Andrés Tuxtla'
$ws.Range('B117').Value = 'Tantoyuca'
$ws.Range('B118').Value = 'Tuxpan'
$ws.Range('B119').Value = 'Veracruz'
$ws.Range('B120').Value = 'Xalapa'
$ws.Range('B121').Value = 'Yanga'
$ws.Range('B122').Value = 'Total'
$ws.Range('A123').Value = 'Total'

# Remove trailing metadata/footer rows (125:129)
$ws.Rows('125:129').Delete()
